$wb = $excel.ActiveWorkbook

# --- Metadata sheet: bump the generation Date ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# --- Elements sheet: swap the two "Mapping" columns (AK <-> AL) ---
$elem = $wb.Worksheets.Item("Elements")

# Swap header text (column titles) - use .Text for reading (the .Value
# getter is not reliable for reads in this host), .Value for writing.
$akHeader = $elem.Cells.Item(1, 37).Text
$alHeader = $elem.Cells.Item(1, 38).Text
$elem.Cells.Item(1, 37).Value = $alHeader
$elem.Cells.Item(1, 38).Value = $akHeader

# Swap the data rows 2-6 between column AK (37) and AL (38)
for ($r = 2; $r -le 6; $r++) {
    $akVal = $elem.Cells.Item($r, 37).Text
    $alVal = $elem.Cells.Item($r, 38).Text
    $elem.Cells.Item($r, 37).Value = $alVal
    $elem.Cells.Item($r, 38).Value = $akVal
}

# Swap the column widths to match the swapped content (AK now holds the
# wide "Spécification métier" text -> ~76.5 chars wide, AL now holds the
# narrower "RIM Mapping" column -> ~25 chars wide).
$elem.Columns.Item(37).ColumnWidth = 75.65
$elem.Columns.Item(38).ColumnWidth = 24.15
